# Auto-generated PowerShell COM-interop script
# Applies the "averaged intensities on spiral schemes" edit:
#  - 3 new orientation schemes are introduced (Spiral-90deg-10rot-5space,
#    Spiral-90deg-15rot-5space, Spiral-90deg-10rot-3space), inserted right
#    after "Gaussian-Quadrature", which itself moves up to directly follow
#    "Ring Perpendicular to TD".
#  - The remaining schemes (NoRotation-tilt60deg ... HexGrid-60degTilt5degRes)
#    shift down to make room, and were re-evaluated (producing tiny
#    floating-point differences vs. the original run for the last 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 4 blank rows before the "NoRotation-tilt60deg" row (row 10). ---
# This shifts the old rows 10-16 down to rows 14-20 (the old "Gaussian-Quadrature"
# row, previously at row 16, becomes a stray duplicate at row 20 - it is removed
# in step 3 below since its data is rewritten into row 10 instead).
$ws.Rows("10:13").Insert()

# Fix the formatting of the 4 newly inserted (blank) rows so column A matches the
# rest of the table (bold + bordered style) instead of Excel's auto-extended style.
$ws.Range("A9:M9").Copy()
$ws.Range("A10:M13").PasteSpecial(-4122)
$excel.CutCopyMode = 0


# --- Step 2: (re)write rows 10-19 completely with the post-edit data. ---
$ws.Range("A10").Value2 = 8
$ws.Range("B10").Value2 = "Gaussian-Quadrature"
$ws.Range("C10").Value2 = 1.283205560933205
$ws.Range("D10").Value2 = 0.5533404005767953
$ws.Range("E10").Value2 = 1.011479690632174
$ws.Range("F10").Value2 = 1.283205560933205
$ws.Range("G10").Value2 = 0.7658860871064203
$ws.Range("H10").Value2 = 1.176224941566672
$ws.Range("I10").Value2 = 1.079696260835131
$ws.Range("J10").Value2 = 0.5533404005767953
$ws.Range("K10").Value2 = 0.7824100456044845
$ws.Range("L10").Value2 = 1.032807803268845
$ws.Range("M10").Value2 = 0.9783054902750662

$ws.Range("A11").Value2 = 9
$ws.Range("B11").Value2 = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value2 = 0.7624447654678239
$ws.Range("D11").Value2 = 1.23810911409098
$ws.Range("E11").Value2 = 1.103482249748861
$ws.Range("F11").Value2 = 0.7624447654678239
$ws.Range("G11").Value2 = 0.8354774746902707
$ws.Range("H11").Value2 = 1.716613816911243
$ws.Range("I11").Value2 = 0.9846972419867169
$ws.Range("J11").Value2 = 1.23810911409098
$ws.Range("K11").Value2 = 1.170795681919921
$ws.Range("L11").Value2 = 0.9666202236938723
$ws.Range("M11").Value2 = 1.106804110482649

$ws.Range("A12").Value2 = 10
$ws.Range("B12").Value2 = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value2 = 0.7634497892880356
$ws.Range("D12").Value2 = 1.241189958470491
$ws.Range("E12").Value2 = 1.102179451484345
$ws.Range("F12").Value2 = 0.7634497892880356
$ws.Range("G12").Value2 = 0.8376267172478166
$ws.Range("H12").Value2 = 1.710581682049699
$ws.Range("I12").Value2 = 0.9836533324204633
$ws.Range("J12").Value2 = 1.241189958470491
$ws.Range("K12").Value2 = 1.171684704977418
$ws.Range("L12").Value2 = 0.9675672471327268
$ws.Range("M12").Value2 = 1.106446821826808

$ws.Range("A13").Value2 = 11
$ws.Range("B13").Value2 = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value2 = 0.7624117444593062
$ws.Range("D13").Value2 = 1.238552621705216
$ws.Range("E13").Value2 = 1.103410334173562
$ws.Range("F13").Value2 = 0.7624117444593062
$ws.Range("G13").Value2 = 0.8360093886495377
$ws.Range("H13").Value2 = 1.715644362247476
$ws.Range("I13").Value2 = 0.9841451338796103
$ws.Range("J13").Value2 = 1.238552621705216
$ws.Range("K13").Value2 = 1.170981477939389
$ws.Range("L13").Value2 = 0.9666966111993476
$ws.Range("M13").Value2 = 1.106695597519118

$ws.Range("A14").Value2 = 12
$ws.Range("B14").Value2 = "NoRotation-tilt60deg"
$ws.Range("C14").Value2 = 0.6059320000000015
$ws.Range("D14").Value2 = 0.7997120000000032
$ws.Range("E14").Value2 = 1.338772
$ws.Range("F14").Value2 = 0.6059320000000015
$ws.Range("G14").Value2 = 0.544116
$ws.Range("H14").Value2 = 2.419648
$ws.Range("I14").Value2 = 1.168527999999998
$ws.Range("J14").Value2 = 0.7997120000000032
$ws.Range("K14").Value2 = 1.069242000000002
$ws.Range("L14").Value2 = 0.8375870000000015
$ws.Range("M14").Value2 = 1.146118

$ws.Range("A15").Value2 = 13
$ws.Range("B15").Value2 = "Rotation-NoTilt"
$ws.Range("C15").Value2 = 0.35
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 1.7072375
$ws.Range("F15").Value2 = 0.35
$ws.Range("G15").Value2 = 0.12
$ws.Range("H15").Value2 = 3.420424999999998
$ws.Range("I15").Value2 = 1.414437499999998
$ws.Range("J15").Value2 = 0
$ws.Range("K15").Value2 = 0.8536187499999999
$ws.Range("L15").Value2 = 0.601809375
$ws.Range("M15").Value2 = 1.168683333333333

$ws.Range("A16").Value2 = 14
$ws.Range("B16").Value2 = "Rotation-60detTilt"
$ws.Range("C16").Value2 = 0.6298105585664012
$ws.Range("D16").Value2 = 0.4086701045760012
$ws.Range("E16").Value2 = 1.4057189262336
$ws.Range("F16").Value2 = 0.6298105585664012
$ws.Range("G16").Value2 = 0.4847611174912008
$ws.Range("H16").Value2 = 2.378300878131199
$ws.Range("I16").Value2 = 1.231107060121596
$ws.Range("J16").Value2 = 0.4086701045760012
$ws.Range("K16").Value2 = 0.9071945154048005
$ws.Range("L16").Value2 = 0.7685025369856008
$ws.Range("M16").Value2 = 1.08972810752

$ws.Range("A17").Value2 = 15
$ws.Range("B17").Value2 = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value2 = 0.980835366581522
$ws.Range("D17").Value2 = 0.9932526162782277
$ws.Range("E17").Value2 = 0.9954731919829907
$ws.Range("F17").Value2 = 0.980835366581522
$ws.Range("G17").Value2 = 0.9879817567710224
$ws.Range("H17").Value2 = 0.9988752707099485
$ws.Range("I17").Value2 = 0.9914000029316566
$ws.Range("J17").Value2 = 0.9932526162782277
$ws.Range("K17").Value2 = 0.9943629041306092
$ws.Range("L17").Value2 = 0.9875991353560656
$ws.Range("M17").Value2 = 0.9913030342092278

$ws.Range("A18").Value2 = 16
$ws.Range("B18").Value2 = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value2 = 1.055120349330215
$ws.Range("D18").Value2 = 1.062051379275293
$ws.Range("E18").Value2 = 0.9434096811110427
$ws.Range("F18").Value2 = 1.055120349330215
$ws.Range("G18").Value2 = 1.012425431463024
$ws.Range("H18").Value2 = 0.938203255348837
$ws.Range("I18").Value2 = 0.9687507810392281
$ws.Range("J18").Value2 = 1.062051379275293
$ws.Range("K18").Value2 = 1.002730530193168
$ws.Range("L18").Value2 = 1.028925439761692
$ws.Range("M18").Value2 = 0.9966601462612732

$ws.Range("A19").Value2 = 17
$ws.Range("B19").Value2 = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value2 = 0.9675179852060344
$ws.Range("D19").Value2 = 1.15609030302424
$ws.Range("E19").Value2 = 0.9443480477397658
$ws.Range("F19").Value2 = 0.9675179852060344
$ws.Range("G19").Value2 = 1.094496109764759
$ws.Range("H19").Value2 = 0.8285166405324971
$ws.Range("I19").Value2 = 0.9431054365710803
$ws.Range("J19").Value2 = 1.15609030302424
$ws.Range("K19").Value2 = 1.050219175382003
$ws.Range("L19").Value2 = 1.008868580294019
$ws.Range("M19").Value2 = 0.9890124204730627

# --- Step 3: remove the stray duplicate row (old Gaussian-Quadrature, pushed to row 20). ---
$ws.Rows("20:20").Delete()

# --- Step 4: tidy up the selection. ---
$ws.Range("A1").Select()
